$wb = $excel.ActiveWorkbook
$cc = $wb.Worksheets.Item("country-code")

# 1) "Euro area" (row 10) is dropped from the dataset entirely.
$cc.Range("A10:C10").Delete()

# 2) Create the new "dropped-country-codes" sheet, positioned between
#    "country-code" and "time-period".
$dropped = $wb.Worksheets.Add($null, $cc)
$dropped.Name = "dropped-country-codes"

# Header row, same as country-code.
$cc.Range("A1:C1").Copy($dropped.Range("A1:C1"))

# 3) Move Hong Kong, Malaysia, Singapore, Taiwan, Thailand out of
#    country-code and into dropped-country-codes (in their original
#    relative order). Row numbers below reflect country-code AFTER the
#    "Euro area" row was already removed above. Cut from the bottom up
#    so row numbers of the still-to-be-moved rows aren't disturbed, but
#    paste each one into its final resting row on the destination sheet.
#    Only the A:C used range is touched (not the whole row) so no stray
#    cells get created out past column C.
$cc.Range("A25:C25").Cut($dropped.Range("A6:C6"))
$cc.Range("A25:C25").Delete()

$cc.Range("A24:C24").Cut($dropped.Range("A5:C5"))
$cc.Range("A24:C24").Delete()

$cc.Range("A19:C19").Cut($dropped.Range("A4:C4"))
$cc.Range("A19:C19").Delete()

$cc.Range("A14:C14").Cut($dropped.Range("A3:C3"))
$cc.Range("A14:C14").Delete()

$cc.Range("A10:C10").Cut($dropped.Range("A2:C2"))
$cc.Range("A10:C10").Delete()

# 4) Selections: country-code becomes the active/selected sheet again,
#    with A10 selected; dropped-country-codes keeps a resting selection
#    at B8; time-period loses its "active" status.
$dropped.Range("B8").Select()
$cc.Activate()
$cc.Range("A10").Select()
